$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (D) and 1h volume-change (E) columns to refreshed values.
# D-column prices are stored as text (not numbers), so force text format before
# assigning to avoid Excel auto-converting them to numeric values (which would
# strip meaningful trailing zeros, e.g. "0.3810" -> 0.381, "1.000" -> 1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.218.89"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.28"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.09"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5367"
$ws.Range("E7").Value = "  +3.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07286"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.17"
$ws.Range("E10").Value = "  +4.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9026"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08197"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.12"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.333"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.83"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008645"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.253.02"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.034"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.096.98"
$ws.Range("E21").Value = "  -42.19%  "
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.490"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.55"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.286"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.32"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.53"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.805"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.751"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09217"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8352"
$ws.Range("E32").Value = "  +4.84%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.216"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.997"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.336"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("E37").Value = "  +3.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5796"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.075"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.347"
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.604"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.85"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1523"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4979"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.635"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.42"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06171"
$ws.Range("E50").Value = "  +3.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.21"
$ws.Range("E51").Value = "  -0.39%  "
